$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, centered, top-aligned, bordered) from H1
# onto the two new header cells I1 and J1, so they match the other headers exactly.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Set the new header labels
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data for the new "I0" (column I) and "IF" (column J) values, one entry per
# data row (rows 2-90).
$i0Values = @(8,7,6,7,8,7,8,7,8,7,7,6,8,8,7,7,8,8,8,8,8,7,7,7,8,7,8,8,8,8,8,8,8,7,8,7,8,8,8,5,8,7,7,7,8,8,6,7,6,7,7,7,6,8,7,7,7,9,7,5,6,6,9,7,7,8,8,7,9,8,7,7,7,6,7,6,7,8,7,9,6,6,6,6,8,6,9,9,3)
$ifValues = @(8,7,7,7,8,7,8,7,8,7,7,6,8,8,8,7,8,8,8,8,8,7,7,8,8,7,8,8,8,8,8,8,8,7,8,7,8,8,8,5,8,7,7,7,8,8,6,7,6,7,7,7,6,8,7,7,7,9,7,6,6,6,9,7,7,8,8,7,9,8,7,7,7,7,7,6,7,8,7,9,6,6,6,6,8,7,9,9,3)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
